# Auto-generated Excel COM-interop script to apply the Bahamut_Profits diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2039.2858  # H40
$ws.Cells.Item(40, 9).Value = 1350  # I40
$ws.Cells.Item(40, 10).Value = 2154.1667  # J40
$ws.Cells.Item(40, 11).Value = 1350  # K40
$ws.Cells.Item(40, 12).Value = 2154.1667  # L40
$ws.Cells.Item(40, 13).Value = -1175  # M40
$ws.Cells.Item(40, 14).Value = -2504.1667  # N40
$ws.Cells.Item(76, 8).Value = 3574503.5  # H76
$ws.Cells.Item(76, 9).Value = 3574503.5  # I76
$ws.Cells.Item(76, 10).Value = 0  # J76
$ws.Cells.Item(76, 11).Value = 3574503.5  # K76
$ws.Cells.Item(76, 12).Value = 0  # L76
$ws.Cells.Item(76, 13).Value = -3574188.5  # M76
$ws.Cells.Item(76, 14).ClearContents()  # N76
$ws.Cells.Item(79, 8).Value = 3574503.5  # H79
$ws.Cells.Item(79, 9).Value = 3574503.5  # I79
$ws.Cells.Item(79, 10).Value = 0  # J79
$ws.Cells.Item(79, 11).Value = 3574503.5  # K79
$ws.Cells.Item(79, 12).Value = 0  # L79
$ws.Cells.Item(79, 13).Value = -3573411.5  # M79
$ws.Cells.Item(79, 14).ClearContents()  # N79
$ws.Cells.Item(98, 8).Value = 2813.8696  # H98
$ws.Cells.Item(98, 9).Value = 2353.6316  # I98
$ws.Cells.Item(98, 10).Value = 5000  # J98
$ws.Cells.Item(98, 11).Value = 2353.6316  # K98
$ws.Cells.Item(98, 12).Value = 5000  # L98
$ws.Cells.Item(98, 13).Value = -855.6316000000002  # M98
$ws.Cells.Item(98, 14).Value = -7996  # N98
$ws.Cells.Item(105, 8).Value = 29000  # H105
$ws.Cells.Item(105, 10).Value = 29000  # J105
$ws.Cells.Item(105, 12).Value = 29000  # L105
$ws.Cells.Item(105, 14).Value = -35988  # N105
$ws.Cells.Item(114, 8).Value = 24756.75  # H114
$ws.Cells.Item(114, 10).Value = 24756.75  # J114
$ws.Cells.Item(114, 12).Value = 24756.75  # L114
$ws.Cells.Item(114, 14).Value = -33434.75  # N114
$ws.Cells.Item(121, 8).Value = 1100  # H121
$ws.Cells.Item(121, 10).Value = 1400  # J121
$ws.Cells.Item(121, 12).Value = 4200  # L121
$ws.Cells.Item(121, 14).Value = -7694  # N121
$ws.Cells.Item(122, 8).Value = 2813.8696  # H122
$ws.Cells.Item(122, 9).Value = 2353.6316  # I122
$ws.Cells.Item(122, 10).Value = 5000  # J122
$ws.Cells.Item(122, 11).Value = 7060.8948  # K122
$ws.Cells.Item(122, 12).Value = 15000  # L122
$ws.Cells.Item(122, 13).Value = -4610.8948  # M122
$ws.Cells.Item(122, 14).Value = -19900  # N122
$ws.Cells.Item(132, 8).Value = 2688.889  # H132
$ws.Cells.Item(132, 9).Value = 3291.0527  # I132
$ws.Cells.Item(132, 10).Value = 1258.75  # J132
$ws.Cells.Item(132, 11).Value = 9873.158100000001  # K132
$ws.Cells.Item(132, 12).Value = 3776.25  # L132
$ws.Cells.Item(132, 13).Value = -7343.158100000001  # M132
$ws.Cells.Item(132, 14).Value = -8836.25  # N132
$ws.Cells.Item(135, 8).Value = 1497.2307  # H135
$ws.Cells.Item(135, 9).Value = 921.7143  # I135
$ws.Cells.Item(135, 10).Value = 3914.4  # J135
$ws.Cells.Item(135, 11).Value = 8295.4287  # K135
$ws.Cells.Item(135, 12).Value = 35229.6  # L135
$ws.Cells.Item(135, 13).Value = -5760.4287  # M135
$ws.Cells.Item(135, 14).Value = -40299.6  # N135
$ws.Cells.Item(137, 8).Value = 832.8421  # H137
$ws.Cells.Item(137, 9).Value = 767.6316  # I137
$ws.Cells.Item(137, 10).Value = 898.0526  # J137
$ws.Cells.Item(137, 11).Value = 2302.8948  # K137
$ws.Cells.Item(137, 12).Value = 2694.1578  # L137
$ws.Cells.Item(137, 13).Value = 247.1052  # M137
$ws.Cells.Item(137, 14).Value = -7794.1578  # N137
$ws.Cells.Item(138, 8).Value = 1735.3  # H138
$ws.Cells.Item(138, 9).Value = 794.2679000000001  # I138
$ws.Cells.Item(138, 10).Value = 2932.9773  # J138
$ws.Cells.Item(138, 11).Value = 2382.8037  # K138
$ws.Cells.Item(138, 12).Value = 8798.9319  # L138
$ws.Cells.Item(138, 13).Value = 2757.1963  # M138
$ws.Cells.Item(138, 14).Value = -19078.9319  # N138
$ws.Cells.Item(141, 8).Value = 2126.3125  # H141
$ws.Cells.Item(141, 9).Value = 783.7273  # I141
$ws.Cells.Item(141, 11).Value = 2351.1819  # K141
$ws.Cells.Item(141, 13).Value = 2828.8181  # M141

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2875.3  # H32
$ws.Cells.Item(32, 9).Value = 3014.6494  # I32
$ws.Cells.Item(32, 10).Value = 2408.7827  # J32
$ws.Cells.Item(32, 11).Value = 3014.6494  # K32
$ws.Cells.Item(32, 12).Value = 2408.7827  # L32
$ws.Cells.Item(32, 13).Value = -2727.6494  # M32
$ws.Cells.Item(32, 14).Value = -2982.7827  # N32
$ws.Cells.Item(61, 8).Value = 1018.8611  # H61
$ws.Cells.Item(61, 9).Value = 876.56525  # I61
$ws.Cells.Item(61, 10).Value = 1270.6154  # J61
$ws.Cells.Item(61, 11).Value = 876.56525  # K61
$ws.Cells.Item(61, 12).Value = 1270.6154  # L61
$ws.Cells.Item(61, 13).Value = -664.56525  # M61
$ws.Cells.Item(61, 14).Value = -1694.6154  # N61
$ws.Cells.Item(74, 8).Value = 1035.9048  # H74
$ws.Cells.Item(74, 9).Value = 966.7879  # I74
$ws.Cells.Item(74, 11).Value = 966.7879  # K74
$ws.Cells.Item(74, 13).Value = -92.78790000000004  # M74
$ws.Cells.Item(77, 8).Value = 1035.9048  # H77
$ws.Cells.Item(77, 9).Value = 966.7879  # I77
$ws.Cells.Item(77, 11).Value = 4833.9395  # K77
$ws.Cells.Item(77, 13).Value = -465.9395000000004  # M77
$ws.Cells.Item(122, 8).Value = 1075  # H122
$ws.Cells.Item(122, 9).Value = 1000  # I122
$ws.Cells.Item(122, 11).Value = 3000  # K122
$ws.Cells.Item(122, 13).Value = -550  # M122
$ws.Cells.Item(132, 8).Value = 1030.8823  # H132
$ws.Cells.Item(132, 9).Value = 942.13794  # I132
$ws.Cells.Item(132, 11).Value = 2826.41382  # K132
$ws.Cells.Item(132, 13).Value = -296.4138199999998  # M132
$ws.Cells.Item(136, 8).Value = 1018.8611  # H136
$ws.Cells.Item(136, 9).Value = 876.56525  # I136
$ws.Cells.Item(136, 10).Value = 1270.6154  # J136
$ws.Cells.Item(136, 11).Value = 2629.69575  # K136
$ws.Cells.Item(136, 12).Value = 3811.8462  # L136
$ws.Cells.Item(136, 13).Value = -79.69574999999986  # M136
$ws.Cells.Item(136, 14).Value = -8911.8462  # N136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 55557410  # H99
$ws.Cells.Item(99, 9).Value = 66668508  # I99
$ws.Cells.Item(99, 10).Value = 1899.6666  # J99
$ws.Cells.Item(99, 11).Value = 66668508  # K99
$ws.Cells.Item(99, 12).Value = 1899.6666  # L99
$ws.Cells.Item(99, 13).Value = -66667010  # M99
$ws.Cells.Item(99, 14).Value = -4895.6666  # N99

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3804.0334  # H31
$ws.Cells.Item(31, 9).Value = 4227.125  # I31
$ws.Cells.Item(31, 10).Value = 3320.5  # J31
$ws.Cells.Item(31, 11).Value = 4227.125  # K31
$ws.Cells.Item(31, 12).Value = 3320.5  # L31
$ws.Cells.Item(31, 13).Value = -3932.125  # M31
$ws.Cells.Item(31, 14).Value = -3910.5  # N31
$ws.Cells.Item(34, 8).Value = 3804.0334  # H34
$ws.Cells.Item(34, 9).Value = 4227.125  # I34
$ws.Cells.Item(34, 10).Value = 3320.5  # J34
$ws.Cells.Item(34, 11).Value = 4227.125  # K34
$ws.Cells.Item(34, 12).Value = 3320.5  # L34
$ws.Cells.Item(34, 13).Value = -4025.125  # M34
$ws.Cells.Item(34, 14).Value = -3724.5  # N34
$ws.Cells.Item(58, 8).Value = 1545.6316  # H58
$ws.Cells.Item(58, 9).Value = 1797.2307  # I58
$ws.Cells.Item(58, 10).Value = 1000.5  # J58
$ws.Cells.Item(58, 11).Value = 1797.2307  # K58
$ws.Cells.Item(58, 12).Value = 1000.5  # L58
$ws.Cells.Item(58, 13).Value = -1594.2307  # M58
$ws.Cells.Item(58, 14).Value = -1406.5  # N58
$ws.Cells.Item(99, 8).Value = 2551.561  # H99
$ws.Cells.Item(99, 9).Value = 2179.1667  # I99
$ws.Cells.Item(99, 10).Value = 3077.2942  # J99
$ws.Cells.Item(99, 11).Value = 2179.1667  # K99
$ws.Cells.Item(99, 12).Value = 3077.2942  # L99
$ws.Cells.Item(99, 13).Value = -681.1667000000002  # M99
$ws.Cells.Item(99, 14).Value = -6073.2942  # N99
$ws.Cells.Item(122, 8).Value = 1200  # H122
$ws.Cells.Item(122, 10).Value = 1200  # J122
$ws.Cells.Item(122, 12).Value = 3600  # L122
$ws.Cells.Item(122, 14).Value = -8500  # N122
$ws.Cells.Item(126, 8).Value = 2551.561  # H126
$ws.Cells.Item(126, 9).Value = 2179.1667  # I126
$ws.Cells.Item(126, 10).Value = 3077.2942  # J126
$ws.Cells.Item(126, 11).Value = 6537.500100000001  # K126
$ws.Cells.Item(126, 12).Value = 9231.882599999999  # L126
$ws.Cells.Item(126, 13).Value = -4067.500100000001  # M126
$ws.Cells.Item(126, 14).Value = -14171.8826  # N126
$ws.Cells.Item(132, 8).Value = 2084.9473  # H132
$ws.Cells.Item(132, 9).Value = 1581.1428  # I132
$ws.Cells.Item(132, 10).Value = 2378.8333  # J132
$ws.Cells.Item(132, 11).Value = 4743.428400000001  # K132
$ws.Cells.Item(132, 12).Value = 7136.499899999999  # L132
$ws.Cells.Item(132, 13).Value = -2213.428400000001  # M132
$ws.Cells.Item(132, 14).Value = -12196.4999  # N132
$ws.Cells.Item(134, 8).Value = 1839.9375  # H134
$ws.Cells.Item(134, 9).Value = 1850.6086  # I134
$ws.Cells.Item(134, 10).Value = 1812.6666  # J134
$ws.Cells.Item(134, 11).Value = 5551.825800000001  # K134
$ws.Cells.Item(134, 12).Value = 5437.9998  # L134
$ws.Cells.Item(134, 13).Value = -3016.825800000001  # M134
$ws.Cells.Item(134, 14).Value = -10507.9998  # N134
$ws.Cells.Item(136, 8).Value = 1545.6316  # H136
$ws.Cells.Item(136, 9).Value = 1797.2307  # I136
$ws.Cells.Item(136, 10).Value = 1000.5  # J136
$ws.Cells.Item(136, 11).Value = 5391.6921  # K136
$ws.Cells.Item(136, 12).Value = 3001.5  # L136
$ws.Cells.Item(136, 13).Value = -2841.6921  # M136
$ws.Cells.Item(136, 14).Value = -8101.5  # N136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 109.71429  # H8
$ws.Cells.Item(8, 9).Value = 109.71429  # I8
$ws.Cells.Item(8, 11).Value = 329.14287  # K8
$ws.Cells.Item(8, 13).Value = -190.14287  # M8
$ws.Cells.Item(68, 8).Value = 681.26666  # H68
$ws.Cells.Item(68, 10).Value = 701.4167  # J68
$ws.Cells.Item(68, 12).Value = 2104.2501  # L68
$ws.Cells.Item(68, 14).Value = -3726.2501  # N68
$ws.Cells.Item(71, 8).Value = 681.26666  # H71
$ws.Cells.Item(71, 10).Value = 701.4167  # J71
$ws.Cells.Item(71, 12).Value = 6312.7503  # L71
$ws.Cells.Item(71, 14).Value = -14424.7503  # N71
$ws.Cells.Item(113, 8).Value = 555.9535  # H113
$ws.Cells.Item(113, 9).Value = 583.0952  # I113
$ws.Cells.Item(113, 10).Value = 530.0454999999999  # J113
$ws.Cells.Item(113, 11).Value = 1749.2856  # K113
$ws.Cells.Item(113, 12).Value = 1590.1365  # L113
$ws.Cells.Item(113, 13).Value = 420.7144000000001  # M113
$ws.Cells.Item(113, 14).Value = -5930.1365  # N113

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 8).Value = 20200  # H51
$ws.Cells.Item(51, 9).Value = 15000  # I51
$ws.Cells.Item(51, 10).Value = 23666.666  # J51
$ws.Cells.Item(51, 11).Value = 15000  # K51
$ws.Cells.Item(51, 12).Value = 23666.666  # L51
$ws.Cells.Item(51, 13).Value = -14491  # M51
$ws.Cells.Item(51, 14).Value = -24684.666  # N51
$ws.Cells.Item(70, 8).Value = 4213.5  # H70
$ws.Cells.Item(70, 9).Value = 3427  # I70
$ws.Cells.Item(70, 11).Value = 3427  # K70
$ws.Cells.Item(70, 13).Value = -3157  # M70
$ws.Cells.Item(73, 8).Value = 4213.5  # H73
$ws.Cells.Item(73, 9).Value = 3427  # I73
$ws.Cells.Item(73, 11).Value = 3427  # K73
$ws.Cells.Item(73, 13).Value = -2491  # M73
$ws.Cells.Item(132, 8).Value = 1897.6595  # H132
$ws.Cells.Item(132, 9).Value = 1799.4849  # I132
$ws.Cells.Item(132, 10).Value = 2129.0715  # J132
$ws.Cells.Item(132, 11).Value = 5398.4547  # K132
$ws.Cells.Item(132, 12).Value = 6387.2145  # L132
$ws.Cells.Item(132, 13).Value = -2868.4547  # M132
$ws.Cells.Item(132, 14).Value = -11447.2145  # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3586930.8  # H7
$ws.Cells.Item(7, 9).Value = 2768.2856  # I7
$ws.Cells.Item(7, 10).Value = 11113672  # J7
$ws.Cells.Item(7, 11).Value = 2768.2856  # K7
$ws.Cells.Item(7, 12).Value = 11113672  # L7
$ws.Cells.Item(7, 13).Value = -2656.2856  # M7
$ws.Cells.Item(7, 14).Value = -11113896  # N7
$ws.Cells.Item(126, 8).Value = 3586930.8  # H126
$ws.Cells.Item(126, 9).Value = 2768.2856  # I126
$ws.Cells.Item(126, 10).Value = 11113672  # J126
$ws.Cells.Item(126, 11).Value = 8304.856800000001  # K126
$ws.Cells.Item(126, 12).Value = 33341016  # L126
$ws.Cells.Item(126, 13).Value = -5834.856800000001  # M126
$ws.Cells.Item(126, 14).Value = -33345956  # N126

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 942.5  # H126
$ws.Cells.Item(126, 9).Value = 750  # I126
$ws.Cells.Item(126, 10).Value = 963.8889  # J126
$ws.Cells.Item(126, 11).Value = 2250  # K126
$ws.Cells.Item(126, 12).Value = 2891.6667  # L126
$ws.Cells.Item(126, 13).Value = 220  # M126
$ws.Cells.Item(126, 14).Value = -7831.6667  # N126
$ws.Cells.Item(132, 8).Value = 1011.3514  # H132
$ws.Cells.Item(132, 9).Value = 993.6667  # I132
$ws.Cells.Item(132, 10).Value = 1087.1428  # J132
$ws.Cells.Item(132, 11).Value = 2981.0001  # K132
$ws.Cells.Item(132, 12).Value = 3261.4284  # L132
$ws.Cells.Item(132, 13).Value = -451.0001000000002  # M132
$ws.Cells.Item(132, 14).Value = -8321.428400000001  # N132
